$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.06854566666666667
$ws.Range("M2").Value = 15.75563966666667
$ws.Range("N2").Value = 47.266919
$ws.Range("O2").Value = 0.3220556913988901
$ws.Range("P2").Value = 0.32205569139889
$ws.Range("Q2").Value = 1.079980824711445
$ws.Range("R2").Value = 9.719827422403
$ws.Range("S2").Value = 0.3220556913988901
$ws.Range("T2").Value = 0.32205569139889

# Row 3
$ws.Range("G3").Value = 0.06854566666666667
$ws.Range("O3").Value = 0.5509544596378365
$ws.Range("P3").Value = 0.5509544596378364
$ws.Range("S3").Value = 0.5509544596378365
$ws.Range("T3").Value = 0.5509544596378364

# Row 4
$ws.Range("G4").Value = 0.06854566666666667
$ws.Range("O4").Value = 0.1269898489632735
$ws.Range("P4").Value = 0.1269898489632735
$ws.Range("Q4").Value = 0.4258474713414445
$ws.Range("S4").Value = 0.1269898489632735
$ws.Range("T4").Value = 0.1269898489632735
